$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Student ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Assignment1"
$ws.Range("D1").Value = "Assignment2"
$ws.Range("E1").Value = "Assignment3"
$ws.Range("F1").Value = "Assignment4"
$ws.Range("G1").Value = "Assignment5"
$ws.Range("H1").Value = "Assignment6"
$ws.Range("I1").Value = "Assignment7"
$ws.Range("J1").Value = "Assignment8"
$ws.Range("K1").Value = "Assignment9"
$ws.Range("L1").Value = "Assignment10"
$ws.Range("M1").Value = "Assignment11"
$ws.Range("N1").Value = "Quiz1"
$ws.Range("O1").Value = "Project1"
$ws.Range("P1").Value = "Pre"
$ws.Range("Q1").Value = "Averaged Score"
$ws.Range("S1").Value = "Note"
$ws.Range("T1").Value = "Email"
$ws.Range("B2").Value = "John Reed"
$ws.Range("T2").Value = "john_reed@mcp.com"
$ws.Range("B3").Value = "Andrew Diaz"
$ws.Range("T3").Value = "diaza48@mcp.com"
$ws.Range("B4").Value = "Matthew Bailey"
$ws.Range("T4").Value = "matthew.bailey32@mcp.com"
$ws.Range("B5").Value = "Ryan Gomez"
$ws.Range("T5").Value = "ryang13@mcp.com"
$ws.Range("B6").Value = "Patricia Watson"
$ws.Range("T6").Value = "pwatson23@mcp.com"
$ws.Range("B7").Value = "Kimberly Murphy"
$ws.Range("T7").Value = "kimberly.murphy@mcp.com"
$ws.Range("B8").Value = "Edward Ruiz"
$ws.Range("T8").Value = "edward.ruiz@mcp.com"
$ws.Range("B9").Value = "Shirley Edwards"
$ws.Range("T9").Value = "shirley_edwards@mcp.com"
$ws.Range("B10").Value = "Catherine Murphy"
$ws.Range("T10").Value = "murphyc@mcp.com"
$ws.Range("B11").Value = "Sandra Miller"
$ws.Range("S11").Value = "Already withdrew from the course"
$ws.Range("T11").Value = "sandra.miller@mcp.com"
$ws.Range("B12").Value = "Carol Carter"
$ws.Range("T12").Value = "carterc23@mcp.com"
$ws.Range("B13").Value = "Michelle Brooks"
$ws.Range("T13").Value = "michelle_brooks26@mcp.com"
$ws.Range("B14").Value = "Steven Morgan"
$ws.Range("T14").Value = "smorgan@mcp.com"
$ws.Range("B15").Value = "Carolyn Alvarez"
$ws.Range("T15").Value = "calvarez@mcp.com"
$ws.Range("B16").Value = "Jennifer Castillo"
$ws.Range("T16").Value = "castilloj@mcp.com"
$ws.Range("B17").Value = "Andrew Moore"
$ws.Range("S17").Value = "Auditing teacher"
$ws.Range("T17").Value = "moorea@mcp.com"
